$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(536).Insert()

$ws.Cells.Item(536, 1).Value = 6
$ws.Cells.Item(536, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(536, 3).Value = "Metropolitana"
$ws.Cells.Item(536, 4).Value = 45041
$ws.Cells.Item(536, 5).Value = 13
$ws.Cells.Item(536, 6).Value = 100112032
$ws.Cells.Item(536, 7).Value = "Zapallo italiano"
$ws.Cells.Item(536, 8).Value = "Sin especificar"
$ws.Cells.Item(536, 9).Value = "Primera"
$ws.Cells.Item(536, 10).Value = 580
$ws.Cells.Item(536, 11).Value = 9000
$ws.Cells.Item(536, 12).Value = 10000
$ws.Cells.Item(536, 13).Value = 9448
$ws.Cells.Item(536, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(536, 15).Value = "Región Metropolitana"
$ws.Cells.Item(536, 16).Value = 189
$ws.Cells.Item(536, 17).Value = 50
$ws.Cells.Item(536, 18).Value = "Hortaliza"
